$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Costos")

function Set-RowRange($sheet, $rangeAddr, $values) {
    $n = $values.Count
    $arr = New-Object 'object[,]' 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $sheet.Range($rangeAddr).Value = $arr
}

$row4Vals = @(7665.5,7857,8053.5,8255,8461.5,8673,8889.5,9112,9339.5,9573,9812.5,13409.5,13744.5,14088,14440.5,14801.5,15171.5,15550.5,15939.5,16338,16746.5,17165,17594,18034,18485,21910,22458,23019.5,23595,24184.5,24789.5,25409,26044.5,26695.5,27363,28047,28748,29466.5,30203.5,30958.5,31732.5,32525.5,33339,34172.5,35026.5,35902.5,36800,37720,38663,39629.5,40620,41635.5,42676.5,43743.5,44837,45958,47107,48284.5,49492,50729,51997.5,53297,54629.5,55995.5,57395.5,58830,60301,61808.5,63353.5,64937.5,66561,68225,69930.5,71679,73471,75307.5,77190.5,79120,81098,83125.5,85203.5,87333.5,89517,91755,94049,96400,98810,101280.5,103812.5,106407.5,109068,111794.5,114589.5,117454,120390.5,123400,126485,129647.5,132888.5,136211,139616,143106.5,146684)
Set-RowRange $ws "K4:DI4" $row4Vals

$row6Vals = @(12757,13076,13403,13738,14081.5,14433.5,14794.5,15164,15543.5,15932,16330,21572,22111.5,22664,23230.5,23811.5,24407,25017,25642.5,26283.5,26940.5,27614,28304.5,29012,29737.5,37568.5,38507.5,39470.5,40457,41468.5,42505,43568,44657,45773.5,46917.5,48090.5,49293,50525.5,51788.5,53083,54410,55770.5,57164.5,58594,60058.5,61560,63099,64676.5,66293.5,67951,69649.5,71391,73175.5,75005,76880,78802,80772.5,82791.5,84861.5,86983,89157.5,91386.5,93671,96013,98413,100873.5,103395.5,105980,108629.5,111345.5,114129,116982.5,119907,122904.5,125977,129126.5,132354.5,135663.5,139055,142531.5,146095,149747,153491,157328,161261.5,165293,169425.5,173661,178002.5,182452.5,187014,191689,196481.5,201393.5,206428,211589,216878.5,222300.5,227858,233554.5,239393.5,245378.5,251512.5)
Set-RowRange $ws "K6:DI6" $row6Vals

$row8Vals = @(1950,1998.5,2048.5,2100,2152.5,2206,2261,2318,2375.5,2435,2496,3191.5,3271.5,3353,3437,3523,3611,3701,3793.5,3888.5,3986,4085.5,4187.5,4292.5,4399.5,5178,5307.5,5440,5576,5715.5,5858.5,6005,6155,6309,6466.5,6628,6794,6963.5,7138,7316,7499,7686.5,7879,8076,8277.5,8484.5,8696.5,8914,9137,9365.5,9599.5,9839.5,10085.5,10337.5,10596,10861,11132.5,11411,11696,11988.5,12288,12595.5,12910.5,13233,13564,13903,14250.5,14607,14972,15346.5,15730,16123,16526.5,16939.5,17363,17797,18242,18698,19165.5,19644.5,20135.5,20639,21155,21684,22226,22781.5,23351,23935,24533.5,25146.5,25775.5,26420,27080.5,27757.5,28451,29162.5,29891.5,30639,31405,32190,32994.5,33819.5,34665)
Set-RowRange $ws "K8:DI8" $row8Vals

$row5Vals = @(0,0,0,0,0,0,7665.5,7857,8053.5,8255,8461.5,8673,8889.5,9112,9339.5,9573,9812.5,13409.5,13744.5,14088,14440.5,14801.5,15171.5,15550.5,15939.5,16338,16746.5,17165,17594,18034,18485,21910,22458,23019.5,23595,24184.5,24789.5,25409,26044.5,26695.5,27363,28047,28748,29466.5,30203.5,30958.5,31732.5,32525.5,33339,34172.5,35026.5,35902.5,36800,37720,38663,39629.5,40620,41635.5,42676.5,43743.5,44837,45958,47107,48284.5,49492,50729,51997.5,53297,54629.5,55995.5,57395.5,58830,60301,61808.5,63353.5,64937.5,66561,68225,69930.5,71679,73471,75307.5,77190.5,79120,81098,83125.5,85203.5,87333.5,89517,91755,94049,96400,98810,101280.5,103812.5,106407.5,109068,111794.5,114589.5,117454,120390.5,123400,126485,129647.5,132888.5,136211,139616,143106.5,146684)
Set-RowRange $ws "E5:DI5" $row5Vals

$row7Vals = @(0,0,0,0,0,0,12757,13076,13403,13738,14081.5,14433.5,14794.5,15164,15543.5,15932,16330,21572,22111.5,22664,23230.5,23811.5,24407,25017,25642.5,26283.5,26940.5,27614,28304.5,29012,29737.5,37568.5,38507.5,39470.5,40457,41468.5,42505,43568,44657,45773.5,46917.5,48090.5,49293,50525.5,51788.5,53083,54410,55770.5,57164.5,58594,60058.5,61560,63099,64676.5,66293.5,67951,69649.5,71391,73175.5,75005,76880,78802,80772.5,82791.5,84861.5,86983,89157.5,91386.5,93671,96013,98413,100873.5,103395.5,105980,108629.5,111345.5,114129,116982.5,119907,122904.5,125977,129126.5,132354.5,135663.5,139055,142531.5,146095,149747,153491,157328,161261.5,165293,169425.5,173661,178002.5,182452.5,187014,191689,196481.5,201393.5,206428,211589,216878.5,222300.5,227858,233554.5,239393.5,245378.5,251512.5)
Set-RowRange $ws "E7:DI7" $row7Vals

$row9Vals = @(0,0,0,0,0,0,1950,1998.5,2048.5,2100,2152.5,2206,2261,2318,2375.5,2435,2496,3191.5,3271.5,3353,3437,3523,3611,3701,3793.5,3888.5,3986,4085.5,4187.5,4292.5,4399.5,5178,5307.5,5440,5576,5715.5,5858.5,6005,6155,6309,6466.5,6628,6794,6963.5,7138,7316,7499,7686.5,7879,8076,8277.5,8484.5,8696.5,8914,9137,9365.5,9599.5,9839.5,10085.5,10337.5,10596,10861,11132.5,11411,11696,11988.5,12288,12595.5,12910.5,13233,13564,13903,14250.5,14607,14972,15346.5,15730,16123,16526.5,16939.5,17363,17797,18242,18698,19165.5,19644.5,20135.5,20639,21155,21684,22226,22781.5,23351,23935,24533.5,25146.5,25775.5,26420,27080.5,27757.5,28451,29162.5,29891.5,30639,31405,32190,32994.5,33819.5,34665)
Set-RowRange $ws "E9:DI9" $row9Vals

$ws.Range("A19:XFD26").Select()